$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Contact Number column (G) holds numeric-looking strings (leading zeros) -
# force text formatting first so they are stored as text, not numbers.
$ws.Range("G8:G12").NumberFormat = "@"

# Row 8 - Malinta National High School (public) entry
$ws.Range("A8").Value = 111222
$ws.Range("B8").Value = "Malinta National High School"
$ws.Range("C8").Value = "public"
$ws.Range("D8").Value = "DCS Valenzuela"
$ws.Range("E8").Value = "Congressional I"
$ws.Range("F8").Value = "MNHS Admin"
$ws.Range("G8").Value = "09123456789"
$ws.Range("H8").Value = "xandrexxenosaquinde@gmail.com"
$ws.Range("I8").Value = "2024-09-03 23:52:39"

# Row 9 - West National High School (private) entry
$ws.Range("A9").Value = 222444
$ws.Range("B9").Value = "West National High School"
$ws.Range("C9").Value = "private"
$ws.Range("D9").Value = "DCS Valenzuela"
$ws.Range("E9").Value = "Congressional I"
$ws.Range("F9").Value = "Lorem Ipsum"
$ws.Range("G9").Value = "09060158736"
$ws.Range("H9").Value = "123@1.com"
$ws.Range("I9").Value = "2024-09-04 01:25:53"

# Row 10 - Dalandanan National High School (public) entry
$ws.Range("A10").Value = 333444
$ws.Range("B10").Value = "Dalandanan National High School"
$ws.Range("C10").Value = "public"
$ws.Range("D10").Value = "DCS Valenzuela"
$ws.Range("E10").Value = "Congressional I"
$ws.Range("F10").Value = 12345
$ws.Range("G10").Value = "09060158736"
$ws.Range("H10").Value = "123@1.com"
$ws.Range("I10").Value = "2024-09-04 01:48:42"

# Row 11 - Maysan National High School (public) entry
$ws.Range("A11").Value = 555666
$ws.Range("B11").Value = "Maysan National High School"
$ws.Range("C11").Value = "public"
$ws.Range("D11").Value = "DCS Valenzuela"
$ws.Range("E11").Value = "Congressional I"
$ws.Range("F11").Value = "Lorem Ipsum A"
$ws.Range("G11").Value = "09060158736"
$ws.Range("H11").Value = "123@1.com"
$ws.Range("I11").Value = "2024-09-04 01:50:22"

# Row 12 (new row) - Example (public) entry
$ws.Range("A12").Value = 999999
$ws.Range("B12").Value = "Example"
$ws.Range("C12").Value = "public"
$ws.Range("D12").Value = "DCS Valenzuela"
$ws.Range("E12").Value = "Congressional I"
$ws.Range("F12").Value = 12345
$ws.Range("G12").Value = "09060158736"
$ws.Range("H12").Value = "123@1.com"
$ws.Range("I12").Value = "2024-09-15 10:39:24"
